$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Build the new Cypher query text for the FilesTab row.
# Note: includes an embedded zero-width space (U+200B) on its own line,
# matching the original authored content exactly.
$zwsp = [char]0x200B
$query = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)<-[*]-(prt)<--(f:file)`n" + `
    "    WHERE c.gender = ""MALE""`n" + `
    "WITH DISTINCT f, prt, c, a, ct`n" + `
    "RETURN`n" + `
    "    COALESCE(f.file_name, '') AS ``File Name``,`n" + `
    "    COALESCE(head(labels(prt)), '') AS ``Association``,`n" + `
    "$zwsp`n" + `
    "    COALESCE(f.file_description, '') AS ``Description``,`n" + `
    "    COALESCE(f.file_format, '') AS ``File Format``,`n" + `
    "    COALESCE(f.file_size, '') AS ``Size``,`n" + `
    "    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,`n" + `
    "    COALESCE(a.arm_id, '') AS ``Arm``,`n" + `
    "    COALESCE(c.case_id, '') AS ``Case ID``"

# Same StatQuery / dbExcel / WebExcel text as row 2 (CasesTab) - reused verbatim.
$statQuery = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)`n" + `
    "    WHERE c.gender = ""MALE""`n" + `
    "OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)`n" + `
    "RETURN `n" + `
    "    COUNT(DISTINCT f) AS number_of_files,`n" + `
    "    COUNT(DISTINCT c.case_id) AS number_of_cases,`n" + `
    "    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials"

$dbExcelFile = "TC01_Trials_Filter_Gender-Male_Neo4jData.xlsx"
$webExcelFile = "TC01_Trials_Filter_Gender-Male_WebData.xlsx"

# Append the new FilesTab row (row 3) under the existing CasesTab row.
$ws.Range("A3").Value = "FilesTab"
$ws.Range("B3").Value = $query
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = $dbExcelFile
$ws.Range("E3").Value = $webExcelFile

# Match the wrap-text style used on B2/C2 for the new B3/C3 cells.
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# Match the row height Excel computed for the wrapped, multi-line query text.
$ws.Rows.Item(3).RowHeight = 188.5

$ws.Range("B3").Select()
